$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need to be forced back to
# Text so they match the source data (which stores every Price/Volume cell as
# text, even when the text looks numeric). Values that are not valid numbers
# (e.g. "31.109.68", which has two dots) are already stored as text by Excel
# automatically and do not need this treatment.

$ws.Range("D2").Value = '31.109.68'
$ws.Range("E2").Value = '  +1.68%  '

$ws.Range("D3").Value = '1.961.58'

$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4902'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2966'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.53%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06885'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.38'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '106.97'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.92%  '

$ws.Range("D12").Value = '1.967.18'
$ws.Range("E12").Value = '  +0.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07805'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.489'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.30%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7031'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '283.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.49%  '

$ws.Range("D17").Value = '31.119.75'
$ws.Range("E17").Value = '  +1.68%  '

$ws.Range("B18").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C18").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D18").Value = '2.248.39'
$ws.Range("E18").Value = '  +1.73%  '

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.07%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007750'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.511'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.519'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.851'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.20%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.206'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1056'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.401'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.579'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.620'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.462'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.22%  '

$ws.Range("E34").Value = '  -2.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7592'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.88%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.175'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.735'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02022'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.704'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.22%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.559'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.23%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.135'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9072'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4488'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.03%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '109.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.202'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.002'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.33%  '

$ws.Range("D48").Value = '1.032.77'
$ws.Range("E48").Value = '  +11.65%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1262'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.356'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.02'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.60%  '
